# Insert a new weekly record row for "Apio" at row 230, shifting all
# subsequent rows (230-265) down to (231-266), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 230; existing rows 230..265 shift to 231..266
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new record's data
$ws.Cells.Item(230, 1).Value = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value = Get-Date -Year 2022 -Month 3 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(230, 5).Value = 9
$ws.Cells.Item(230, 6).Value = 100112017
$ws.Cells.Item(230, 7).Value = "Apio"
$ws.Cells.Item(230, 8).Value = "Americana (o)"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 250
$ws.Cells.Item(230, 11).Value = 10000
$ws.Cells.Item(230, 12).Value = 12000
$ws.Cells.Item(230, 13).Value = 10800
$ws.Cells.Item(230, 14).Value = "`$/docena de matas"
$ws.Cells.Item(230, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(230, 16).Value = 1800
$ws.Cells.Item(230, 17).Value = 6
$ws.Cells.Item(230, 18).Value = "Hortaliza"
